$wb = $excel.ActiveWorkbook

# A scratch cell (way outside the used range of every sheet we touch)
# used to coerce a date-looking string ("2024-05-07" etc.) into a
# literal text value before it is pasted into the real target cell.
# Using NumberFormat "@" directly on the target cell and then trying
# to reset its format back to General always ends up registering a
# brand new style entry, so instead we prepare the text on a
# throw-away cell, copy *only the value* over with PasteSpecial, and
# wipe the scratch cell again so it leaves no trace behind.
function Set-TextCellStyled($ws, $range, [string]$text) {
    $scratch = $ws.Range("ZZ100")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $range.Style = "Normal"
    $scratch.Copy() | Out-Null
    $range.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues) | Out-Null
    $excel.CutCopyMode = 0
    $scratch.Clear() | Out-Null
}

# Same idea, but for cells that must end up with no explicit style at
# all (plain default formatting), matching row 21 in the diff.
function Set-TextCellPlain($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

function Set-DateTimeCellStyled($range, [double]$value) {
    $range.Style = "Normal"
    $range.Value2 = $value
    $range.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

function Set-NumberCellStyled($range, $value) {
    $range.Style = "Normal"
    $range.Value = $value
}

# ---------------------------------------------------------------
# Sheet: AMSIN
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("AMSIN")

# Fix row 16: correct the run time recorded in B16 and normalize the
# formatting of the rest of the row.
$ws.Range("A16").Style = "Normal"
$ws.Range("C16:G16").Style = "Normal"
$ws.Range("B16").Value2 = 45202.49348515047

# New row 17
Set-TextCellStyled $ws $ws.Range("A17") "2024-05-07"
Set-DateTimeCellStyled $ws.Range("B17") 45419.75788628472
Set-TextCellStyled $ws $ws.Range("C17") "191aaccpp"
Set-NumberCellStyled $ws.Range("D17") 62
Set-NumberCellStyled $ws.Range("E17") 30
Set-NumberCellStyled $ws.Range("F17") 32
Set-NumberCellStyled $ws.Range("G17") 1

# New row 18
Set-TextCellStyled $ws $ws.Range("A18") "2024-05-08"
Set-DateTimeCellStyled $ws.Range("B18") 45420.44737079861
Set-TextCellStyled $ws $ws.Range("C18") "191accpp"
Set-NumberCellStyled $ws.Range("D18") 62
Set-NumberCellStyled $ws.Range("E18") 30
Set-NumberCellStyled $ws.Range("F18") 32
Set-NumberCellStyled $ws.Range("G18") 0.83

# New row 19
Set-TextCellStyled $ws $ws.Range("A19") "2024-05-08"
Set-DateTimeCellStyled $ws.Range("B19") 45420.47590496528
Set-TextCellStyled $ws $ws.Range("C19") "191prorp"
Set-NumberCellStyled $ws.Range("D19") 34
Set-NumberCellStyled $ws.Range("E19") 30
Set-NumberCellStyled $ws.Range("F19") 4
Set-NumberCellStyled $ws.Range("G19") 1.24

# New row 20
Set-TextCellStyled $ws $ws.Range("A20") "2024-05-08"
Set-DateTimeCellStyled $ws.Range("B20") 45420.47833754629
Set-TextCellStyled $ws $ws.Range("C20") "191kkacp"
Set-NumberCellStyled $ws.Range("D20") 34
Set-NumberCellStyled $ws.Range("E20") 34
Set-NumberCellStyled $ws.Range("F20") 0
Set-NumberCellStyled $ws.Range("G20") 0.91

# New row 21 (A/C/D/E/F/G left without an explicit style, like the
# source diff, but B21 keeps the same run-time/date style as B17-B20)
Set-TextCellPlain $ws.Range("A21") "2024-05-08"
Set-DateTimeCellStyled $ws.Range("B21") 45420.4847673906
Set-TextCellPlain $ws.Range("C21") "191finalacp"
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = 35
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 1.12

# ---------------------------------------------------------------
# Sheet: AMS
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("AMS")

# New row 15
Set-TextCellStyled $ws2 $ws2.Range("A15") "2024-05-08"
Set-DateTimeCellStyled $ws2.Range("B15") 45420.45077725694
Set-TextCellStyled $ws2 $ws2.Range("C15") "191vinacp"
Set-NumberCellStyled $ws2.Range("D15") 34
Set-NumberCellStyled $ws2.Range("E15") 30
Set-NumberCellStyled $ws2.Range("F15") 4
Set-NumberCellStyled $ws2.Range("G15") 0.89
